$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel
# (pure numeric-looking strings) must be pre-formatted as Text so the
# stored value stays an exact string match, like the source data feed.
$textCells = @(
    "D5",
    "D7",
    "D8",
    "D10",
    "D11",
    "D12",
    "D13",
    "D15",
    "D16",
    "D19",
    "D23",
    "D24",
    "D25",
    "D26",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D34",
    "D35",
    "D36",
    "D49",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D46",
    "D47",
    "D48",
    "D50",
    "D51",
)
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# --- Simple price/volume updates (rows with unchanged coin identity) ---
$ws.Range("D2").Value = "30.420.55"
$ws.Range("E2").Value = "  +0.64%  "
$ws.Range("D3").Value = "1.869.06"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("D5").Value = "246.01"
$ws.Range("E5").Value = "  +1.79%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").Value = "0.4737"
$ws.Range("E7").Value = "  +0.71%  "
$ws.Range("D8").Value = "0.2917"
$ws.Range("E8").Value = "  +2.30%  "
$ws.Range("E9").Value = "  +0.40%  "
$ws.Range("D10").Value = "22.07"
$ws.Range("E10").Value = "  +6.29%  "
$ws.Range("D11").Value = "0.07718"
$ws.Range("E11").Value = "  +0.18%  "
$ws.Range("D12").Value = "97.58"
$ws.Range("E12").Value = "  +2.54%  "
$ws.Range("D13").Value = "0.7395"
$ws.Range("E13").Value = "  +8.19%  "
$ws.Range("D14").Value = "1.872.20"
$ws.Range("E14").Value = "  +0.46%  "
$ws.Range("D15").Value = "5.130"
$ws.Range("E15").Value = "  +0.81%  "
$ws.Range("D16").Value = "273.63"
$ws.Range("E16").Value = "  +1.74%  "
$ws.Range("D17").Value = "30.408.71"
$ws.Range("E17").Value = "  +0.67%  "
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").Value = "0.000007537"
$ws.Range("E19").Value = "  +0.22%  "
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").Value = "2.117.82"
$ws.Range("E21").Value = "  +0.53%  "
$ws.Range("D23").Value = "5.222"
$ws.Range("E23").Value = "  +0.25%  "
$ws.Range("D24").Value = "6.172"
$ws.Range("E24").Value = "  +0.83%  "
$ws.Range("D25").Value = "9.306"
$ws.Range("E25").Value = "  -0.27%  "
$ws.Range("D26").Value = "164.18"
$ws.Range("E26").Value = "  -1.01%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("E28").Value = "  +1.80%  "
$ws.Range("D29").Value = "0.09983"
$ws.Range("E29").Value = "  +1.69%  "
$ws.Range("D30").Value = "1.366"
$ws.Range("E30").Value = "  -0.59%  "
$ws.Range("D31").Value = "1.500"
$ws.Range("E31").Value = "  -0.23%  "
$ws.Range("D32").Value = "4.302"
$ws.Range("E32").Value = "  +1.31%  "
$ws.Range("D33").Value = "4.136"
$ws.Range("E33").Value = "  +3.66%  "
$ws.Range("D34").Value = "0.04840"
$ws.Range("E34").Value = "  +3.21%  "
$ws.Range("D35").Value = "1.122"
$ws.Range("E35").Value = "  +0.86%  "
$ws.Range("D36").Value = "0.6979"
$ws.Range("E36").Value = "  +1.95%  "
$ws.Range("D49").Value = "6.999"
$ws.Range("E49").Value = "  +0.75%  "

# --- Rows 37-51: ranking reshuffled, full row content replaced ---
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "2.714"
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.01859"
$ws.Range("E38").Value = "  +1.16%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "2.743"
$ws.Range("E39").Value = "  +0.52%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "6.309"
$ws.Range("E40").Value = "  -0.46%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "73.17"
$ws.Range("E41").Value = "  +3.87%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "1.967"
$ws.Range("E42").Value = "  +4.25%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "0.4198"
$ws.Range("E43").Value = "  +3.27%  "
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").Value = "0.8342"
$ws.Range("E45").Value = "  -0.53%  "
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "102.08"
$ws.Range("E46").Value = "  +0.23%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "9.222"
$ws.Range("E47").Value = "  -0.08%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "931.25"
$ws.Range("E48").Value = "  +0.97%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "35.43"
$ws.Range("E50").Value = "  +2.73%  "
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").Value = "0.3884"
$ws.Range("E51").Value = "  +3.04%  "
